$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired data for rows 2..18 (header in row 1 stays unchanged)
$data = @(
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Scotty Pippen Jr.", "PG,SG", "Memphis Grizzlies"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Jalen Johnson", "PF", "Atlanta Hawks"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Mark Williams", "C", "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-obsolete last row (old row 19, Darius Garland) so the
# sheet shrinks from 19 rows to 18 rows.
$ws.Rows(19).Delete()
